# Update "想去人数" (interest count) figures in the 广州-漫展信息 workbook.
# Sheet "展览" (Exhibitions)
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 2751
$wsExpo.Range("F4").Value = 1083
$wsExpo.Range("F5").Value = 19960
$wsExpo.Range("F7").Value = 2318
$wsExpo.Range("F10").Value = 449
$wsExpo.Range("F12").Value = 250
$wsExpo.Range("F19").Value = 212

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 295
$wsShow.Range("F16").Value = 100

# Sheet "全部类型" (All types - aggregated list)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 2751
$wsAll.Range("F9").Value = 1083
$wsAll.Range("F10").Value = 19960
$wsAll.Range("F15").Value = 295
$wsAll.Range("F16").Value = 2318
$wsAll.Range("F20").Value = 449
$wsAll.Range("F22").Value = 250
$wsAll.Range("F36").Value = 212
$wsAll.Range("F37").Value = 100
$wsAll.Range("F38").Value = 100
